$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 37, shifting existing rows 37-58 down to 38-59
$ws.Rows.Item(37).Insert()

# Populate the new row 37 with the new data record
$ws.Cells.Item(37, 1).Value = 11
$ws.Cells.Item(37, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(37, 3).Value = "Bíobío"
$ws.Cells.Item(37, 4).Value = 45126
$ws.Cells.Item(37, 4).NumberFormat = $ws.Cells.Item(38, 4).NumberFormat
$ws.Cells.Item(37, 5).Value = 8
$ws.Cells.Item(37, 6).Value = 100114007
$ws.Cells.Item(37, 7).Value = "Jengibre"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 30
$ws.Cells.Item(37, 11).Value = 16000
$ws.Cells.Item(37, 12).Value = 16000
$ws.Cells.Item(37, 13).Value = 16000
$ws.Cells.Item(37, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(37, 15).Value = "Perú"
$ws.Cells.Item(37, 16).Value = 1231
$ws.Cells.Item(37, 17).Value = 13
$ws.Cells.Item(37, 18).Value = "Hortaliza"
